# "add check in make changes"
# The scheduling algorithm ("make changes") got an extra validation check,
# which produced a different (longer, re-shuffled) shift assignment and
# updated per-worker shift totals.

$wb = $excel.ActiveWorkbook

$wsScrewed = $wb.Worksheets.Item("screwed")
$wsShifts  = $wb.Worksheets.Item("shiftsPerWorker")

# ---------------------------------------------------------------------
# Sheet "screwed": re-generated worker-per-shift assignment, now with
# two extra shifts (rows 13 and 14 / index 12 and 13).
# ---------------------------------------------------------------------
$assignments = @(
    "adir",
    "yoni",
    "rotem",
    "tair",
    "stav",
    "adir",
    "asaf",
    "stav",
    "asaf",
    "rotem",
    "adir",
    "stav",
    "rotem"
)

for ($i = 0; $i -lt $assignments.Length; $i++) {
    $row = $i + 2
    $wsScrewed.Cells.Item($row, 1).Value = $i + 1
    $wsScrewed.Cells.Item($row, 2).Value = $assignments[$i]
}

# ---------------------------------------------------------------------
# Sheet "shiftsPerWorker": updated totals per worker.
# ---------------------------------------------------------------------
$wsShifts.Range("B2").Value = 7
$wsShifts.Range("B3").Value = 3
$wsShifts.Range("B4").Value = 4
$wsShifts.Range("B5").Value = 4
$wsShifts.Range("B6").Value = 2
$wsShifts.Range("B7").Value = 4

# ---------------------------------------------------------------------
# Workbook-wide default font: Calibri -> Arial
# ---------------------------------------------------------------------
$wb.Styles.Item("Normal").Font.Name = "Arial"

# ---------------------------------------------------------------------
# Leave a selection behind on "screwed" and make "shiftsPerWorker" the
# active sheet/tab, as in the saved workbook.
# ---------------------------------------------------------------------
[void]$wsScrewed.Range("E20").Select()
[void]$wsShifts.Activate()
